$wb = $excel.ActiveWorkbook

# --- Sheet: missing_values ---
$ws1 = $wb.Worksheets.Item("missing_values")
$ws1.Range("B38").Value = 102
$ws1.Range("C38").Value = 1.2749999999999999
$ws1.Range("B40").Value = 7834
$ws1.Range("C40").Value = 97.924999999999997

# --- Sheet: labor_incmon_imp_stochastic_reg ---
$ws3 = $wb.Worksheets.Item("labor_incmon_imp_stochastic_reg")
$ws3.Range("A3").Value = 2007357.2352610789
$ws3.Range("B3").Value = 250197.34375
$ws3.Range("C3").Value = 455767.578125
$ws3.Range("F3").Value = 4003157
$ws3.Range("G3").Value = 2079280.2081727833
$ws3.Range("H3").Value = 270213.125
$ws3.Range("I3").Value = 500000
$ws3.Range("J3").Value = 1053030.25
$ws3.Range("K3").Value = 2653803.75
$ws3.Range("L3").Value = 4203804

# --- Sheet: labor_jubpenimp_stochastic_reg ---
$ws4 = $wb.Worksheets.Item("labor_jubpenimp_stochastic_reg")
$ws4.Range("A3").Value = 1518900.9281340449
$ws4.Range("C3").Value = 250197.34375
$ws4.Range("E3").Value = 2454001
$ws4.Range("F3").Value = 2554080
$ws4.Range("G3").Value = 1517309.3726285503
$ws4.Range("I3").Value = 250197.34375
$ws4.Range("K3").Value = 2454001
$ws4.Range("L3").Value = 2554080

# --- Sheet: nonlabor_imp_stochastic_reg ---
$ws5 = $wb.Worksheets.Item("nonlabor_imp_stochastic_reg")
$ws5.Range("A3").Value = 434814.74970929569
$ws5.Range("C3").Value = 200157.875
$ws5.Range("E3").Value = 484000
$ws5.Range("F3").Value = 840515.125
$ws5.Range("G3").Value = 434861.08810818463
$ws5.Range("H3").Value = 120094.71875
$ws5.Range("I3").Value = 200157.875
$ws5.Range("K3").Value = 486968.09375
$ws5.Range("L3").Value = 840515.125

# --- Sheet: labor_beneimp_stochastic_reg ---
$ws6 = $wb.Worksheets.Item("labor_beneimp_stochastic_reg")
$ws6.Range("A3").Value = 866184.52546126244
$ws6.Range("B3").Value = 55043.4140625
$ws6.Range("C3").Value = 188879.8125
$ws6.Range("D3").Value = 300236.8125
$ws6.Range("E3").Value = 856565.625
$ws6.Range("G3").Value = 894069.06448190357
$ws6.Range("I3").Value = 195153.921875
$ws6.Range("J3").Value = 310810.9375
$ws6.Range("K3").Value = 870270.625
$ws6.Range("L3").Value = 2000000
